# Gandhi.docx "reset" edit
#
# 1. The first paragraph currently reads:
#      "This is a Microsoft word document.  (This is a change – Version for branch alternate)"
#    It must be trimmed back to just:
#      "This is a Microsoft word document."
#    i.e. the two trailing spaces plus the red "(This is a change – Version
#    for branch alternate)" run(s) are removed, leaving the single plain run.
#
# 2. The document currently ends with a stray empty paragraph
#    (an empty <w:p/>) right before the final section break. That
#    paragraph must be removed so the section break immediately follows
#    the preceding (NormalWeb) paragraph.

$d = $word.ActiveDocument

# --- Step 1: trim paragraph 1 back down to the plain sentence ---------
$keepText = "This is a Microsoft word document."
$p1 = $d.Paragraphs(1).Range
$p1Text = $p1.Text

if ($p1Text.Length -gt ($keepText.Length + 1)) {
    $delStart = $p1.Start + $keepText.Length
    $delEnd = $p1.End - 1   # stop before the paragraph mark itself
    $extraRange = $d.Range($delStart, $delEnd)
    $extraRange.Delete()
}

# --- Step 2: remove the trailing empty paragraph -----------------------
# Only do this when the very last paragraph is empty (just a paragraph
# mark) and the one before it is empty too, i.e. a redundant blank
# paragraph was tacked on at the end of the body after the intentional
# trailing (NormalWeb) blank paragraph.
$count = $d.Paragraphs.Count
if ($count -ge 2) {
    $lastPara = $d.Paragraphs($count)
    $prevPara = $d.Paragraphs($count - 1)

    if (($lastPara.Range.Text -eq [char]13) -and ($prevPara.Range.Text -eq [char]13)) {
        $killRange = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
        $killRange.Delete()
    }
}
